$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Establish correct per-column cell formatting (style A=1,B=2,C=3) for all
# rows we are about to touch, by copying format from a known-good donor row (9).
# This keeps existing correct styles intact and fixes styles for brand-new cells.
$donor = $ws.Range("A9:C9")
$donor.Copy()
foreach ($r in 10,13,14,15,16,17,18,19,20,21,22,23) {
    $ws.Range("A$r`:C$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Row 10: Objetivos (A unchanged) -> new B/C text
$ws.Range("B10").Value = "Desenvolver conhecimentos de forma a tornar o aluno capaz de interpretar corretamente o desenho técnico, conhecer as metodologias e ferramentas utilizadas na indústria, dando subsídios para que possa executar, interagir e modificar desenhos e projetos ao longo de sua vida profissional."
$ws.Range("C10").Value = "Desenvolver conhecimentos de forma a tornar o aluno capaz de interpretar corretamente o desenho técnico, conhecer as metodologias e ferramentas utilizadas na indústria, dando subsídios para que possa executar, interagir e modificar desenhos e projetos ao longo de sua vida profissional."

# Row 13: teacher name 1 (A empty)
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"

# Row 14: teacher name 2 (A empty)
$ws.Range("A14").ClearContents()
$ws.Range("B14").Value = "5817692 - Katia Cristiane Gandolpho Candioto"
$ws.Range("C14").Value = "5817692 - Katia Cristiane Gandolpho Candioto"

# Row 15: Programa resumido / short syllabus (PT)
$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = "Contexto do desenho técnico na indústria, principais ferramentas e técnicas utilizadas em desenhos para elaboração de projetos. Introdução ao desenho assistido por computador (CAD)."
$ws.Range("C15").Value = "Contexto do desenho técnico na indústria, principais ferramentas e técnicas utilizadas em desenhos para elaboração de projetos. Introdução ao desenho assistido por computador (CAD)."

# Row 16: Short syllabus (EN)
$ws.Range("A16").Value = "Short syllabus:"
$ws.Range("B16").Value = "Context of the technical drawing in the industry, main tools and techniques used in drawings for the elaboration of projects. Introduction to computer-aided design (CAD)."
$ws.Range("C16").Value = "Context of the technical drawing in the industry, main tools and techniques used in drawings for the elaboration of projects. Introduction to computer-aided design (CAD)."

# Row 17: Programa (PT)
$ws.Range("A17").Value = "Programa:"
$ws.Range("B17").Value = "Normas do desenho técnico. Terminologia técnica e materiais para desenho. Representação em perspectiva. Projeto ortogonal. Dimensionamento e escala. Corte e secção. Vista Auxiliar e detalhes. Tolerâncias geométricas. Representação de elementos de máquinas. Utilização de software para desenho técnico. Desenho assistido por computador em três dimensões (Modelagem de Sólidos). Desenho assistido por computador em duas dimensões."
$ws.Range("C17").Value = "Normas do desenho técnico. Terminologia técnica e materiais para desenho. Representação em perspectiva. Projeto ortogonal. Dimensionamento e escala. Corte e secção. Vista Auxiliar e detalhes. Tolerâncias geométricas. Representação de elementos de máquinas. Utilização de software para desenho técnico. Desenho assistido por computador em três dimensões (Modelagem de Sólidos). Desenho assistido por computador em duas dimensões."

# Row 18: Syllabus (EN)
$ws.Range("A18").Value = "Syllabus:"
$ws.Range("B18").Value = "Technical drawing standards. Technical terminology and materials for drawing. Perspective representation. Orthogonal design. Scaling and scaling. Cut and section. Auxiliary view and details. Geometric tolerances. Representation of machine elements. Use of software for technical design. Computer-aided design in three dimensions (Solid Modeling). Computer-aided design in two dimensions."
$ws.Range("C18").Value = "Technical drawing standards. Technical terminology and materials for drawing. Perspective representation. Orthogonal design. Scaling and scaling. Cut and section. Auxiliary view and details. Geometric tolerances. Representation of machine elements. Use of software for technical design. Computer-aided design in three dimensions (Solid Modeling). Computer-aided design in two dimensions."

# Row 19: Avaliacao (header only, B/C empty)
$ws.Range("A19").Value = "Avaliação:"
$ws.Range("B19").ClearContents()
$ws.Range("C19").ClearContents()

# Row 20: Metodo
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "Aulas expositivas, trabalhos e aulas práticas. Aulas com softwares para desenho técnico."
$ws.Range("C20").Value = "Aulas expositivas, trabalhos e aulas práticas. Aulas com softwares para desenho técnico."

# Row 21: Criterio
$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "Média aritmética das notas de atividades em aula e extra aula."
$ws.Range("C21").Value = "Média aritmética das notas de atividades em aula e extra aula."

# Row 22: Norma de recuperacao (new row)
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "Devido às características práticas da disciplina, não será oferecida recuperação"
$ws.Range("C22").Value = "Devido às características práticas da disciplina, não será oferecida recuperação"

# Row 23: Bibliografia (new row)
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "FRENCH, T. E.; VIERCK, C. J. Desenho Técnico e Tecnologia Gráfica, Editora Globo, 1999.`nGIESECKE, F. E. Comunicação Gráfica Moderna, Editora Bookman, 2002.`nRIBEIRO, A. C.; PERES, M. P.; IZIDORO, N. Curso de Desenho Técnico e AutoCAD, Pearson, 2013.`nSILVA, A.; RIBEIRO, C. T.; DIAS, J.; SOUSA, L. Desenho Técnico Moderno, LTC, 2013.`nCRUZ, M. D. Catia V5r20 - Modelagem, Montagem e Detalhamento, ERICA, 2010.`nLIMA, C.C. Estudo dirigido de AutoCAD 2015. ÉRICA, 2015.`nLEAKE, J. Manual de Desenho técnico para engenharia, LTC, 2010.`nFISCHER, U; GOMERINGER, R; HEINZLER, M; ET AL. Manual de Tecnologia Metal Mecânica, Blucher, 2011.`nPROVENZA, F. Desenhista de Máquinas . Editora Protec, 1991.`nPROVENZA, F. Projetista de Máquinas . Editora Protec, 1991."
$ws.Range("C23").Value = "FRENCH, T. E.; VIERCK, C. J. Desenho Técnico e Tecnologia Gráfica, Editora Globo, 1999.`nGIESECKE, F. E. Comunicação Gráfica Moderna, Editora Bookman, 2002.`nRIBEIRO, A. C.; PERES, M. P.; IZIDORO, N. Curso de Desenho Técnico e AutoCAD, Pearson, 2013.`nSILVA, A.; RIBEIRO, C. T.; DIAS, J.; SOUSA, L. Desenho Técnico Moderno, LTC, 2013.`nCRUZ, M. D. Catia V5r20 - Modelagem, Montagem e Detalhamento, ERICA, 2010.`nLIMA, C.C. Estudo dirigido de AutoCAD 2015. ÉRICA, 2015.`nLEAKE, J. Manual de Desenho técnico para engenharia, LTC, 2010.`nFISCHER, U; GOMERINGER, R; HEINZLER, M; ET AL. Manual de Tecnologia Metal Mecânica, Blucher, 2011.`nPROVENZA, F. Desenhista de Máquinas . Editora Protec, 1991.`nPROVENZA, F. Projetista de Máquinas . Editora Protec, 1991."

# --- Row heights to match target layout
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120

